$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.761.51"
$ws.Range("E2").Value = "  -0.40%  "

$ws.Range("D3").Value = "2.555.15"
$ws.Range("E3").Value = "  +0.33%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").Value = "'303.46"
$ws.Range("E5").Value = "  +1.81%  "

$ws.Range("D6").Value = "'98.68"
$ws.Range("E6").Value = "  +7.27%  "

$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.548"
$ws.Range("E9").Value = "  -0.33%  "

$ws.Range("D10").Value = "'36.40"
$ws.Range("E10").Value = "  +1.84%  "

$ws.Range("D11").Value = "'0.0809"
$ws.Range("E11").Value = "  +0.35%  "

$ws.Range("E12").Value = "  +7.64%  "

$ws.Range("D13").Value = "'7.69"
$ws.Range("E13").Value = "  +0.31%  "

$ws.Range("D14").Value = "2.535.73"
$ws.Range("E14").Value = "  -0.32%  "

$ws.Range("D15").Value = "'0.882"
$ws.Range("E15").Value = "  +1.76%  "

$ws.Range("D16").Value = "'14.80"
$ws.Range("E16").Value = "  +4.59%  "

$ws.Range("D17").Value = "42.817.41"
$ws.Range("E17").Value = "  -0.28%  "

$ws.Range("D18").Value = "'13.21"
$ws.Range("E18").Value = "  +5.38%  "

$ws.Range("E19").Value = "  +0.79%  "

$ws.Range("D20").Value = "'6.60"
$ws.Range("E20").Value = "  -0.75%  "

$ws.Range("D21").Value = "'71.67"
$ws.Range("E21").Value = "  -0.52%  "

$ws.Range("D22").Value = "'254.91"
$ws.Range("E22").Value = "  -2.22%  "

$ws.Range("D23").Value = "'2.96"
$ws.Range("E23").Value = "  +1.95%  "

$ws.Range("E24").Value = "  -2.06%  "

$ws.Range("D25").Value = "'27.66"
$ws.Range("E25").Value = "  -6.13%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").Value = "'10.07"
$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D28").Value = "'37.94"
$ws.Range("E28").Value = "  +4.04%  "

$ws.Range("E29").Value = "  -1.35%  "

$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("D31").Value = "'156.76"
$ws.Range("E31").Value = "  +2.68%  "

$ws.Range("E32").Value = "  +0.57%  "

$ws.Range("E33").Value = "  +0.62%  "

$ws.Range("D34").Value = "'0.0809"
$ws.Range("E34").Value = "  +2.14%  "

$ws.Range("E35").Value = "  -2.67%  "

$ws.Range("D36").Value = "'26.48"
$ws.Range("E36").Value = "  +10.20%  "

$ws.Range("D37").Value = "'18.63"
$ws.Range("E37").Value = "  +11.20%  "

$ws.Range("E38").Value = "  +0.90%  "

$ws.Range("E39").Value = "  -0.04%  "

$ws.Range("D40").Value = "'2.09"
$ws.Range("E40").Value = "  +32.38%  "

$ws.Range("D41").Value = "'3.41"
$ws.Range("E41").Value = "  -1.32%  "

$ws.Range("E42").Value = "  +0.86%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0304"
$ws.Range("E43").Value = "  -2.50%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.063.34"
$ws.Range("E44").Value = "  -0.50%  "

$ws.Range("D46").Value = "'87.75"
$ws.Range("E46").Value = "  +2.99%  "

$ws.Range("E47").Value = "  +6.30%  "

$ws.Range("D48").Value = "2.801.66"
$ws.Range("E48").Value = "  +0.27%  "

$ws.Range("D49").Value = "'75.08"
$ws.Range("E49").Value = "  +8.22%  "

$ws.Range("D50").Value = "'103.22"
$ws.Range("E50").Value = "  -0.90%  "

$ws.Range("E51").Value = "  +1.82%  "
